$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $style = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $style
}

Set-TextValue $ws.Range("D2") "304.20"
Set-TextValue $ws.Range("E2") "5.89%"
Set-TextValue $ws.Range("D3") "32.51"
Set-TextValue $ws.Range("E3") "11.61%"
Set-TextValue $ws.Range("D4") "5.291"
Set-TextValue $ws.Range("E4") "1.49%"
Set-TextValue $ws.Range("D5") "0.07496"
Set-TextValue $ws.Range("E5") "7.20%"
Set-TextValue $ws.Range("D6") "7.845"
Set-TextValue $ws.Range("E6") "5.84%"
Set-TextValue $ws.Range("D7") "3.799"
Set-TextValue $ws.Range("E7") "6.82%"
Set-TextValue $ws.Range("D8") "1.518"
Set-TextValue $ws.Range("E8") "8.30%"
Set-TextValue $ws.Range("D9") "0.9173"
Set-TextValue $ws.Range("E9") "2.58%"
Set-TextValue $ws.Range("D10") "0.01774"
Set-TextValue $ws.Range("E10") "2,637.99%"
Set-TextValue $ws.Range("D11") "0.1701"
Set-TextValue $ws.Range("E11") "5.66%"
Set-TextValue $ws.Range("D12") "0.07846"
Set-TextValue $ws.Range("E12") "3.61%"
Set-TextValue $ws.Range("D13") "0.08070"
Set-TextValue $ws.Range("E13") "5.20%"
Set-TextValue $ws.Range("D14") "0.03014"
Set-TextValue $ws.Range("E14") "3.13%"
Set-TextValue $ws.Range("D15") "0.09913"
Set-TextValue $ws.Range("E15") "10.09%"
Set-TextValue $ws.Range("D16") "0.001495"
Set-TextValue $ws.Range("E16") "-5.73%"
Set-TextValue $ws.Range("D17") "0.04616"
Set-TextValue $ws.Range("E17") "2.08%"
Set-TextValue $ws.Range("D18") "0.006172"
Set-TextValue $ws.Range("E18") "-4.12%"
Set-TextValue $ws.Range("D19") "3.469"
Set-TextValue $ws.Range("E19") "0.27%"
Set-TextValue $ws.Range("D20") "2.229"
Set-TextValue $ws.Range("E20") "-0.01%"
Set-TextValue $ws.Range("E21") "2.91%"
Set-TextValue $ws.Range("D22") "0.1335"
Set-TextValue $ws.Range("E22") "0.40%"
Set-TextValue $ws.Range("D23") "4.480"
Set-TextValue $ws.Range("E23") "11.76%"
Set-TextValue $ws.Range("D24") "0.1619"
Set-TextValue $ws.Range("E24") "4.38%"
Set-TextValue $ws.Range("D25") "0.001218"
Set-TextValue $ws.Range("E25") "0.99%"
Set-TextValue $ws.Range("D26") "0.004459"
Set-TextValue $ws.Range("E26") "5.18%"
Set-TextValue $ws.Range("E27") "19.85%"
Set-TextValue $ws.Range("E28") "7.23%"
Set-TextValue $ws.Range("D40") "0.04545"
Set-TextValue $ws.Range("E40") "4.84%"
Set-TextValue $ws.Range("D41") "0.007184"
Set-TextValue $ws.Range("E41") "3.88%"
Set-TextValue $ws.Range("E42") "8.27%"
Set-TextValue $ws.Range("D43") "0.002205"
Set-TextValue $ws.Range("E43") "6.42%"
Set-TextValue $ws.Range("D44") "0.01275"
Set-TextValue $ws.Range("E44") "11.17%"
Set-TextValue $ws.Range("D45") "0.00006198"
Set-TextValue $ws.Range("E45") "6.41%"
Set-TextValue $ws.Range("D46") "1.872"
Set-TextValue $ws.Range("E46") "-2.95%"
Set-TextValue $ws.Range("E47") "-0.51%"
